# Actualización automática 2025-11-19 10:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO -------------------------------------------------
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("D27").Value = 2136.96
$wsVentasGrupo.Range("L27").Value = 1108.6
$wsVentasGrupo.Range("M49").Value = 4804.82
$wsVentasGrupo.Range("D60").Value = "6 de 58"
$wsVentasGrupo.Range("L60").Value = "8 de 58"

# --- Sheet: VENTA MENSUAL ----------------------------------------------------
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F27").Value = 3314.82
$wsVentaMensual.Range("F49").Value = 4804.82
$wsVentaMensual.Range("F60").Value = 24282.69

# --- Sheet: CUMPLIMIENTO MENSUAL --------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D3").Value = 7973.56
$wsCumplimiento.Range("E3").Value = 6851.849999999999
$wsCumplimiento.Range("F3").Value = 0.5378306569599087

$wsCumplimiento.Range("D11").Value = 5020.76
$wsCumplimiento.Range("E11").Value = 11127.24
$wsCumplimiento.Range("F11").Value = 0.3109214763438197

$wsCumplimiento.Range("D12").Value = 11135.11
$wsCumplimiento.Range("E12").Value = 39171.89
$wsCumplimiento.Range("F12").Value = 0.2213431530403324

$wsCumplimiento.Range("D14").Value = 25854.82
$wsCumplimiento.Range("E14").Value = 72007.06766749099
$wsCumplimiento.Range("F14").Value = 0.2641970292648339
